# Controle_Financeiro.xlsx edit:
# - Swap descriptions/payment method on GASTOS between the "Shopping" single
#   purchase and the "Presente" installment purchase (D2/G2 and D3/D4/D5).
# - Make RESUMO the active/selected sheet (was ENTRADAS).

$wb = $excel.ActiveWorkbook

$gastos = $wb.Worksheets.Item("GASTOS")
$resumo = $wb.Worksheets.Item("RESUMO")

# Update GASTOS data: the single "Shopping" purchase becomes "Presente" and
# its payment method changes from Pix to Débito; the 3-installment purchase
# that used to be "Presente (x/3)" becomes "Shopping (x/3)".
$gastos.Range("D2").Value = "Presente"
$gastos.Range("G2").Value = "Débito"
$gastos.Range("D3").Value = "Shopping (1/3)"
$gastos.Range("D4").Value = "Shopping (2/3)"
$gastos.Range("D5").Value = "Shopping (3/3)"

# Switch the active sheet to RESUMO (was ENTRADAS).
$resumo.Activate()

$wb.Save()
